# Generate Report for Handoff
# Update the "Latest Handback DateTime" (column D, row 5) on the zh-cn and
# de-de language sheets with new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-25 05:59:13"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-25 05:59:24"
